$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 164  # H6: was 130.16667
$ws.Cells.Item(6, 9).Value = 52.333332  # I6: was 56.4
$ws.Cells.Item(6, 11).Value = 156.999996  # K6: was 169.2
$ws.Cells.Item(6, 13).Value = -44.99999600000001  # M6: was -57.19999999999999

# Sheet ALC, Row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 9063.656000000001  # H62: was 9110.593999999999
$ws.Cells.Item(62, 10).Value = 11331  # J62: was 11438.286
$ws.Cells.Item(62, 12).Value = 11331  # L62: was 11438.286
$ws.Cells.Item(62, 14).Value = -12579  # N62: was -12686.286

# Sheet ALC, Row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 4227.4  # H64: was 2949.8
$ws.Cells.Item(64, 9).Value = 4227.4  # I64: was 2949.8
$ws.Cells.Item(64, 11).Value = 4227.4  # K64: was 2949.8
$ws.Cells.Item(64, 13).Value = -3979.4  # M64: was -2701.8

# Sheet ALC, Row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 9063.656000000001  # H65: was 9110.593999999999
$ws.Cells.Item(65, 10).Value = 11331  # J65: was 11438.286
$ws.Cells.Item(65, 12).Value = 56655  # L65: was 57191.43
$ws.Cells.Item(65, 14).Value = -62895  # N65: was -63431.43

# Sheet ALC, Row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(67, 8).Value = 4227.4  # H67: was 2949.8
$ws.Cells.Item(67, 9).Value = 4227.4  # I67: was 2949.8
$ws.Cells.Item(67, 11).Value = 4227.4  # K67: was 2949.8
$ws.Cells.Item(67, 13).Value = -3369.4  # M67: was -2091.8

# Sheet ALC, Row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(103, 8).Value = 957.2857  # H103: was 1000.5
$ws.Cells.Item(103, 10).Value = 1166  # J103: was 1400
$ws.Cells.Item(103, 12).Value = 3498  # L103: was 4200
$ws.Cells.Item(103, 14).Value = -4670  # N103: was -5372

# Sheet ALC, Row 115
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(115, 8).Value = 786  # H115: was 948.8570999999999
$ws.Cells.Item(115, 9).Value = 633.6667  # I115: was 773.6667
$ws.Cells.Item(115, 10).Value = 1700  # J115: was 2000
$ws.Cells.Item(115, 11).Value = 1901.0001  # K115: was 2321.0001
$ws.Cells.Item(115, 12).Value = 5100  # L115: was 6000
$ws.Cells.Item(115, 13).Value = -334.0001  # M115: was -754.0001000000002
$ws.Cells.Item(115, 14).Value = -8234  # N115: was -9134

# Sheet ALC, Row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 39451.21  # H132: was 41437.168
$ws.Cells.Item(132, 9).Value = 44596.574  # I132: was 47196.727
$ws.Cells.Item(132, 11).Value = 133789.722  # K132: was 141590.181
$ws.Cells.Item(132, 13).Value = -131259.722  # M132: was -139060.181

# Sheet ARM, Row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6103341.5  # H32: was 6255942.5
$ws.Cells.Item(32, 9).Value = 6497688.5  # I32: was 6670979
$ws.Cells.Item(32, 11).Value = 6497688.5  # K32: was 6670979
$ws.Cells.Item(32, 13).Value = -6497401.5  # M32: was -6670692

# Sheet ARM, Row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 10363.069  # H61: was 10994.407
$ws.Cells.Item(61, 9).Value = 9120.666999999999  # I61: was 9887.053
$ws.Cells.Item(61, 11).Value = 9120.666999999999  # K61: was 9887.053
$ws.Cells.Item(61, 13).Value = -8908.666999999999  # M61: was -9675.053

# Sheet ARM, Row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 2704.2083  # H74: was 2761.4856
$ws.Cells.Item(74, 9).Value = 2382.0784  # I74: was 2450.7551
$ws.Cells.Item(74, 11).Value = 2382.0784  # K74: was 2450.7551
$ws.Cells.Item(74, 13).Value = -1508.0784  # M74: was -1576.7551

# Sheet ARM, Row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 2704.2083  # H77: was 2761.4856
$ws.Cells.Item(77, 9).Value = 2382.0784  # I77: was 2450.7551
$ws.Cells.Item(77, 11).Value = 11910.392  # K77: was 12253.7755
$ws.Cells.Item(77, 13).Value = -7542.392  # M77: was -7885.7755

# Sheet ARM, Row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 3189.2693  # H102: was 3433.3044
$ws.Cells.Item(102, 9).Value = 2606.75  # I102: was 2834.1177
$ws.Cells.Item(102, 11).Value = 2606.75  # K102: was 2834.1177
$ws.Cells.Item(102, 13).Value = -984.75  # M102: was -1212.1177

# Sheet ARM, Row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 10363.069  # H136: was 10994.407
$ws.Cells.Item(136, 9).Value = 9120.666999999999  # I136: was 9887.053
$ws.Cells.Item(136, 11).Value = 27362.001  # K136: was 29661.159
$ws.Cells.Item(136, 13).Value = -24812.001  # M136: was -27111.159

# Sheet BSM, Row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 438.90475  # H94: was 453.9
$ws.Cells.Item(94, 9).Value = 453.35  # I94: was 469.89474
$ws.Cells.Item(94, 11).Value = 453.35  # K94: was 469.89474
$ws.Cells.Item(94, 13).Value = -2.350000000000023  # M94: was -18.89474000000001

# Sheet CRP, Row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 651.1429000000001  # H16: was 735.4
$ws.Cells.Item(16, 9).Value = 630.1667  # I16: was 725
$ws.Cells.Item(16, 11).Value = 630.1667  # K16: was 725
$ws.Cells.Item(16, 13).Value = -343.1667  # M16: was -438

# Sheet CRP, Row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5592.0527  # H31: was 5670.482
$ws.Cells.Item(31, 9).Value = 2075.3  # I31: was 2105.4827
$ws.Cells.Item(31, 11).Value = 2075.3  # K31: was 2105.4827
$ws.Cells.Item(31, 13).Value = -1780.3  # M31: was -1810.4827

# Sheet CRP, Row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 5592.0527  # H34: was 5670.482
$ws.Cells.Item(34, 9).Value = 2075.3  # I34: was 2105.4827
$ws.Cells.Item(34, 11).Value = 2075.3  # K34: was 2105.4827
$ws.Cells.Item(34, 13).Value = -1873.3  # M34: was -1903.4827

# Sheet CRP, Row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2620.3845  # H58: was 2497.1428
$ws.Cells.Item(58, 9).Value = 1924.091  # I58: was 1838.3334
$ws.Cells.Item(58, 11).Value = 1924.091  # K58: was 1838.3334
$ws.Cells.Item(58, 13).Value = -1721.091  # M58: was -1635.3334

# Sheet CRP, Row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 651.1429000000001  # H113: was 735.4
$ws.Cells.Item(113, 9).Value = 630.1667  # I113: was 725
$ws.Cells.Item(113, 11).Value = 630.1667  # K113: was 725
$ws.Cells.Item(113, 13).Value = 1539.8333  # M113: was 1445

# Sheet CRP, Row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 743.375  # H122: was 769.5714
$ws.Cells.Item(122, 9).Value = 763.8570999999999  # I122: was 797.8333
$ws.Cells.Item(122, 11).Value = 2291.5713  # K122: was 2393.4999
$ws.Cells.Item(122, 13).Value = 158.4287000000004  # M122: was 56.5001000000002

# Sheet CRP, Row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 13891360  # H132: was 12502394
$ws.Cells.Item(132, 9).Value = 2342.0715  # I132: was 2261.8125
$ws.Cells.Item(132, 11).Value = 7026.2145  # K132: was 6785.4375
$ws.Cells.Item(132, 13).Value = -4496.2145  # M132: was -4255.4375

# Sheet CRP, Row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 2091.7932  # H134: was 2134.4285
$ws.Cells.Item(134, 9).Value = 1987.7693  # I134: was 2031.36
$ws.Cells.Item(134, 11).Value = 5963.3079  # K134: was 6094.08
$ws.Cells.Item(134, 13).Value = -3428.3079  # M134: was -3559.08

# Sheet CRP, Row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 2620.3845  # H136: was 2497.1428
$ws.Cells.Item(136, 9).Value = 1924.091  # I136: was 1838.3334
$ws.Cells.Item(136, 11).Value = 5772.272999999999  # K136: was 5515.0002
$ws.Cells.Item(136, 13).Value = -3222.272999999999  # M136: was -2965.0002

# Sheet CUL, Row 37
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 98316.164  # H37: was 97128
$ws.Cells.Item(37, 10).Value = 98316.164  # J37: was 97128
$ws.Cells.Item(37, 12).Value = 294948.492  # L37: was 291384
$ws.Cells.Item(37, 14).Value = -295172.492  # N37: was -291608

# Sheet CUL, Row 109
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(109, 8).Value = 3270.6843  # H109: was 3154.111
$ws.Cells.Item(109, 9).Value = 1805.5834  # I109: was 1819
$ws.Cells.Item(109, 10).Value = 5782.2856  # J109: was 5824.3335
$ws.Cells.Item(109, 11).Value = 5416.7502  # K109: was 5457
$ws.Cells.Item(109, 12).Value = 17346.8568  # L109: was 17473.0005
$ws.Cells.Item(109, 13).Value = -4376.7502  # M109: was -4417
$ws.Cells.Item(109, 14).Value = -19426.8568  # N109: was -19553.0005

# Sheet CUL, Row 114
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(114, 8).Value = 1875.5385  # H114: was 1804.1333
$ws.Cells.Item(114, 9).Value = 540.2143  # I114: was 607.625
$ws.Cells.Item(114, 10).Value = 3433.4167  # J114: was 3171.5715
$ws.Cells.Item(114, 11).Value = 1620.6429  # K114: was 1822.875
$ws.Cells.Item(114, 12).Value = 10300.2501  # L114: was 9514.7145
$ws.Cells.Item(114, 13).Value = 1633.3571  # M114: was 1431.125
$ws.Cells.Item(114, 14).Value = -16808.2501  # N114: was -16022.7145

# Sheet CUL, Row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1664.973  # H131: was 1695.9736
$ws.Cells.Item(131, 9).Value = 778.5  # I131: was 763.8
$ws.Cells.Item(131, 10).Value = 2340.3809  # J131: was 2303.913
$ws.Cells.Item(131, 11).Value = 2335.5  # K131: was 2291.4
$ws.Cells.Item(131, 12).Value = 7021.1427  # L131: was 6911.739
$ws.Cells.Item(131, 13).Value = 2704.5  # M131: was 2748.6
$ws.Cells.Item(131, 14).Value = -17101.1427  # N131: was -16991.739

# Sheet CUL, Row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 2550.5  # H137: was 2400.3635
$ws.Cells.Item(137, 9).Value = 2954.6667  # I137: was 2661
$ws.Cells.Item(137, 11).Value = 8864.000100000001  # K137: was 7983
$ws.Cells.Item(137, 13).Value = -3764.000100000001  # M137: was -2883

# Sheet GSM, Row 27
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(27, 8).Value = 0  # H27: was 5000
$ws.Cells.Item(27, 9).Value = 0  # I27: was 5000
$ws.Cells.Item(27, 11).Value = 0  # K27: was 5000
$ws.Cells.Item(27, 13).ClearContents()  # M27: was -4834

# Sheet GSM, Row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1939.8387  # H97: was 2059.5862
$ws.Cells.Item(97, 9).Value = 344.4  # I97: was 355
$ws.Cells.Item(97, 10).Value = 3435.5625  # J97: was 3650.5334
$ws.Cells.Item(97, 11).Value = 344.4  # K97: was 355
$ws.Cells.Item(97, 12).Value = 3435.5625  # L97: was 3650.5334
$ws.Cells.Item(97, 13).Value = 151.6  # M97: was 141
$ws.Cells.Item(97, 14).Value = -4427.5625  # N97: was -4642.5334

# Sheet GSM, Row 109
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(109, 8).Value = 50000  # H109: was 0
$ws.Cells.Item(109, 10).Value = 50000  # J109: was 0
$ws.Cells.Item(109, 12).Value = 50000  # L109: was 0
$ws.Cells.Item(109, 14).Value = -52080  # N109: was None

# Sheet GSM, Row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3353.5518  # H122: was 3416.0356
$ws.Cells.Item(122, 9).Value = 1730.64  # I122: was 1735.9166
$ws.Cells.Item(122, 11).Value = 5191.92  # K122: was 5207.7498
$ws.Cells.Item(122, 13).Value = -2741.92  # M122: was -2757.7498

# Sheet GSM, Row 133
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(133, 8).Value = 88487.25  # H133: was 88980
$ws.Cells.Item(133, 10).Value = 88487.25  # J133: was 88980
$ws.Cells.Item(133, 12).Value = 88487.25  # L133: was 88980
$ws.Cells.Item(133, 14).Value = -98607.25  # N133: was -99100

# Sheet GSM, Row 139
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(139, 8).Value = 188750  # H139: was 185833.33
$ws.Cells.Item(139, 10).Value = 188750  # J139: was 185833.33
$ws.Cells.Item(139, 12).Value = 188750  # L139: was 185833.33
$ws.Cells.Item(139, 14).Value = -199030  # N139: was -196113.33

# Sheet LTW, Row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 12911.25  # H16: was 14702.714
$ws.Cells.Item(16, 9).Value = 12911.25  # I16: was 14702.714
$ws.Cells.Item(16, 11).Value = 12911.25  # K16: was 14702.714
$ws.Cells.Item(16, 13).Value = -12741.25  # M16: was -14532.714

# Sheet LTW, Row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 754480.8  # H132: was 788709.4399999999
$ws.Cells.Item(132, 9).Value = 1050039.4  # I132: was 1117690.2
$ws.Cells.Item(132, 11).Value = 3150118.2  # K132: was 3353070.6
$ws.Cells.Item(132, 13).Value = -3147588.2  # M132: was -3350540.6

# Sheet WVR, Row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 598971.25  # H132: was 609470.8
$ws.Cells.Item(132, 9).Value = 888870.7  # I132: was 912249
$ws.Cells.Item(132, 11).Value = 2666612.1  # K132: was 2736747
$ws.Cells.Item(132, 13).Value = -2664082.1  # M132: was -2734217

# Sheet WVR, Row 138
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(138, 8).Value = 94996.336  # H138: was 94995.5
$ws.Cells.Item(138, 10).Value = 94996.336  # J138: was 94995.5
$ws.Cells.Item(138, 12).Value = 94996.336  # L138: was 94995.5
$ws.Cells.Item(138, 14).Value = -105276.336  # N138: was -105275.5
